# "Updated stuff before logo"
#
# Slide 7 ("Donate to support") carries four donation-tier callouts next to
# badge icons. Three of the four "Title for those who ..." captions are
# reworded to the "donate <amount> ... development of the app." phrasing:
#
#   idx=18 placeholder: "...contribute over Rs.100" -> "...donate Rs.100 amount
#                        towards the development of the app." and a new
#                        trailing empty paragraph is appended (matching the
#                        idx=21 placeholder's existing shape).
#   idx=21 placeholder: "...donate any amount..."    -> "...donate Rs.20 amount
#                        towards the development of the app."
#   idx=22 placeholder: "...contribute over Rs.200"  -> "...donate Rs.500 amount
#                        towards the development of the app."

$rupee = [char]0x20B9

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }

    $tr = $shp.TextFrame.TextRange
    $txt = $tr.Text.TrimEnd("`r")

    if ($txt -eq "Title for those who contribute over " + $rupee + "100") {
        # Reword the single run in place (keeps its rPr / the rest of the
        # shape untouched), then append a brand-new trailing empty paragraph.
        $tr.Runs(1, 1).Text = "Title for those who donate " + $rupee + "100 amount towards the development of the app."
        [void]$shp.TextFrame.TextRange.InsertAfter("`rX")
        $lastPara = $shp.TextFrame.TextRange.Paragraphs(2, 1)
        [void]$lastPara.Delete()
    }
    elseif ($txt -eq "Title for those who donate any amount towards the development of the app.") {
        # Only the wording of the first paragraph changes; the existing
        # trailing empty paragraph must be left exactly as-is.
        $tr.Runs(1, 1).Text = "Title for those who donate " + $rupee + "20 amount towards the development of the app."
    }
    elseif ($txt -eq "Title for those who contribute over " + $rupee + "200") {
        $tr.Runs(1, 1).Text = "Title for those who donate " + $rupee + "500 amount towards the development of the app."
    }
}
